$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task "Minimap & New track" (row 14) is being split into two completed
# tasks: "New track" (row 14) and "Minimap" (row 15). Insert a new row 15,
# shifting "Car Selection & Customize" and everything below it down by one
# (dimension grows from A1:B21 to A1:B22).
$ws.Rows.Item(15).Insert(-4121)

# Copy the formatting (borders/fonts) of the row we just split from into the
# newly inserted row so both halves look identical.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(15).RowHeight = 26.25

# Copy the "DONE" look (green fill, centered) from B2 onto the three newly
# completed status cells (B13, B14, B15).
$ws.Range("B2").Copy()
$ws.Range("B13:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 13: Cutscenes is now complete.
$ws.Range("B13").Value = "DONE"

# Row 14: New track (first half of the split task) - complete.
$ws.Range("A14").Value = "New track"
$ws.Range("B14").Value = "DONE"

# Row 15: Minimap (second half of the split task) - complete.
$ws.Range("A15").Value = "Minimap"
$ws.Range("B15").Value = "DONE"

# Update the view to match where the user was working.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("B13:B14").Select()

$wb.Save()
